# Update the "Förändrad" (Changed) date column (C) for all data rows.
# All values move from the Excel serial date 45202 (2023-10-03)
# to 45203 (2023-10-04), mirroring the automatic daily refresh recorded
# in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}
